$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7821
$ws.Range("J51").Value = 2676.8462
$ws.Range("L51").Value = 2676.8462
$ws.Range("N51").Value = -3644.8462

$ws.Range("H98").Value = 422.88235
$ws.Range("I98").Value = 406.7143
$ws.Range("J98").Value = 498.33334
$ws.Range("K98").Value = 406.7143
$ws.Range("L98").Value = 498.33334
$ws.Range("M98").Value = 1091.2857
$ws.Range("N98").Value = -3494.33334

$ws.Range("H116").Value = 3677.6667
$ws.Range("J116").Value = 4516.5
$ws.Range("L116").Value = 4516.5
$ws.Range("N116").Value = -11400.5

$ws.Range("H122").Value = 422.88235
$ws.Range("I122").Value = 406.7143
$ws.Range("J122").Value = 498.33334
$ws.Range("K122").Value = 1220.1429
$ws.Range("L122").Value = 1495.00002
$ws.Range("M122").Value = 1229.8571
$ws.Range("N122").Value = -6395.000019999999

$ws.Range("H132").Value = 7359851.5
$ws.Range("I132").Value = 7819530
$ws.Range("K132").Value = 23458590
$ws.Range("M132").Value = -23456060

$ws.Range("H136").Value = 58571.43
$ws.Range("J136").Value = 58571.43
$ws.Range("L136").Value = 58571.43
$ws.Range("N136").Value = -68771.42999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6388.296
$ws.Range("I32").Value = 5037.5737
$ws.Range("K32").Value = 5037.5737
$ws.Range("M32").Value = -4750.5737

$ws.Range("H61").Value = 2586
$ws.Range("I61").Value = 1965
$ws.Range("J61").Value = 3000
$ws.Range("K61").Value = 1965
$ws.Range("L61").Value = 3000
$ws.Range("M61").Value = -1753
$ws.Range("N61").Value = -3424

$ws.Range("H74").Value = 873.66
$ws.Range("I74").Value = 830.6667
$ws.Range("J74").Value = 1099.375
$ws.Range("K74").Value = 830.6667
$ws.Range("L74").Value = 1099.375
$ws.Range("M74").Value = 43.33330000000001
$ws.Range("N74").Value = -2847.375

$ws.Range("H77").Value = 873.66
$ws.Range("I77").Value = 830.6667
$ws.Range("J77").Value = 1099.375
$ws.Range("K77").Value = 4153.3335
$ws.Range("L77").Value = 5496.875
$ws.Range("M77").Value = 214.6665000000003
$ws.Range("N77").Value = -14232.875

$ws.Range("H97").Value = 42850.082
$ws.Range("I97").Value = 63256.812
$ws.Range("J97").Value = 2036.625
$ws.Range("K97").Value = 63256.812
$ws.Range("L97").Value = 2036.625
$ws.Range("M97").Value = -62760.812
$ws.Range("N97").Value = -3028.625

$ws.Range("H110").Value = 83426230
$ws.Range("I110").Value = 125138470
$ws.Range("K110").Value = 125138470
$ws.Range("M110").Value = -125136425

$ws.Range("H132").Value = 2447.7715
$ws.Range("I132").Value = 2157.6553
$ws.Range("J132").Value = 3850
$ws.Range("K132").Value = 6472.965899999999
$ws.Range("L132").Value = 11550
$ws.Range("M132").Value = -3942.965899999999
$ws.Range("N132").Value = -16610

$ws.Range("H136").Value = 2586
$ws.Range("I136").Value = 1965
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 5895
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3345
$ws.Range("N136").Value = -14100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 517.9259
$ws.Range("I94").Value = 500
$ws.Range("J94").Value = 537.2308
$ws.Range("K94").Value = 500
$ws.Range("L94").Value = 537.2308
$ws.Range("M94").Value = -49
$ws.Range("N94").Value = -1439.2308

$ws.Range("H107").Value = 90910180
$ws.Range("I107").Value = 250000370
$ws.Range("J107").Value = 1492.2858
$ws.Range("K107").Value = 250000370
$ws.Range("L107").Value = 1492.2858
$ws.Range("M107").Value = -249998450
$ws.Range("N107").Value = -5332.2858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2461.0925
$ws.Range("I31").Value = 1606.5454
$ws.Range("J31").Value = 3048.5938
$ws.Range("K31").Value = 1606.5454
$ws.Range("L31").Value = 3048.5938
$ws.Range("M31").Value = -1311.5454
$ws.Range("N31").Value = -3638.5938

$ws.Range("H34").Value = 2461.0925
$ws.Range("I34").Value = 1606.5454
$ws.Range("J34").Value = 3048.5938
$ws.Range("K34").Value = 1606.5454
$ws.Range("L34").Value = 3048.5938
$ws.Range("M34").Value = -1404.5454
$ws.Range("N34").Value = -3452.5938

$ws.Range("H86").Value = 3477.7144
$ws.Range("I86").Value = 3200
$ws.Range("K86").Value = 3200
$ws.Range("M86").Value = -2077

$ws.Range("H89").Value = 3477.7144
$ws.Range("I89").Value = 3200
$ws.Range("K89").Value = 16000
$ws.Range("M89").Value = -10384

$ws.Range("H99").Value = 12100.272
$ws.Range("I99").Value = 4558.4
$ws.Range("J99").Value = 18385.166
$ws.Range("K99").Value = 4558.4
$ws.Range("L99").Value = 18385.166
$ws.Range("M99").Value = -3060.4
$ws.Range("N99").Value = -21381.166

$ws.Range("H107").Value = 1260.2
$ws.Range("I107").Value = 1458.8
$ws.Range("J107").Value = 1061.6
$ws.Range("K107").Value = 1458.8
$ws.Range("L107").Value = 1061.6
$ws.Range("M107").Value = 461.2
$ws.Range("N107").Value = -4901.6

$ws.Range("H126").Value = 12100.272
$ws.Range("I126").Value = 4558.4
$ws.Range("J126").Value = 18385.166
$ws.Range("K126").Value = 13675.2
$ws.Range("L126").Value = 55155.49800000001
$ws.Range("M126").Value = -11205.2
$ws.Range("N126").Value = -60095.49800000001

$ws.Range("H132").Value = 5768.8945
$ws.Range("I132").Value = 8789
$ws.Range("J132").Value = 3572.4546
$ws.Range("K132").Value = 26367
$ws.Range("L132").Value = 10717.3638
$ws.Range("M132").Value = -23837
$ws.Range("N132").Value = -15777.3638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 11500

$ws.Range("H90").Value = 11500

$ws.Range("H131").Value = 769.02
$ws.Range("I131").Value = 250
$ws.Range("J131").Value = 779.61224
$ws.Range("K131").Value = 750
$ws.Range("L131").Value = 2338.83672
$ws.Range("M131").Value = 4290
$ws.Range("N131").Value = -12418.83672

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 5885266
$ws.Range("J126").Value = 8405209
$ws.Range("L126").Value = 25215627
$ws.Range("N126").Value = -25220567

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H88").Value = 6840
$ws.Range("I88").Value = 6840
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 6840
$ws.Range("L88").Value = 0
$ws.Range("M88").Value = -6412
$ws.Range("N88").Value = ""

$ws.Range("H91").Value = 6840
$ws.Range("I91").Value = 6840
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 6840
$ws.Range("L91").Value = 0
$ws.Range("M91").Value = -5358
$ws.Range("N91").Value = ""

$ws.Range("H122").Value = 2791.1667
$ws.Range("I122").Value = 2724.25
$ws.Range("J122").Value = 2925
$ws.Range("K122").Value = 8172.75
$ws.Range("L122").Value = 8775
$ws.Range("M122").Value = -5722.75
$ws.Range("N122").Value = -13675

$ws.Range("H136").Value = 1798.8572
$ws.Range("I136").Value = 1610.3529
$ws.Range("J136").Value = 2600
$ws.Range("K136").Value = 4831.0587
$ws.Range("L136").Value = 7800
$ws.Range("M136").Value = -2281.0587
$ws.Range("N136").Value = -12900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2254.9
$ws.Range("I122").Value = 1478.1111
$ws.Range("J122").Value = 2890.4546
$ws.Range("K122").Value = 4434.3333
$ws.Range("L122").Value = 8671.363799999999
$ws.Range("M122").Value = -1984.3333
$ws.Range("N122").Value = -13571.3638

$ws.Range("H126").Value = 1603.9
$ws.Range("I126").Value = 1317.375
$ws.Range("J126").Value = 2750
$ws.Range("K126").Value = 3952.125
$ws.Range("L126").Value = 8250
$ws.Range("M126").Value = -1482.125
$ws.Range("N126").Value = -13190

$ws.Range("H132").Value = 1743.2833
$ws.Range("I132").Value = 1715.5745
$ws.Range("J132").Value = 1843.4615
$ws.Range("K132").Value = 5146.7235
$ws.Range("L132").Value = 5530.3845
$ws.Range("M132").Value = -2616.7235
$ws.Range("N132").Value = -10590.3845

$ws.Range("H136").Value = 1500.303
$ws.Range("I136").Value = 596.9167
$ws.Range("J136").Value = 3909.3333
$ws.Range("K136").Value = 1790.7501
$ws.Range("L136").Value = 11727.9999
$ws.Range("M136").Value = 759.2499
$ws.Range("N136").Value = -16827.9999
